# Day 6: Min Heap Construction, DFS, Caesar Cipher Encryptor
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Dynamic Programming sheet: add column H note + fill E7/F7
# ---------------------------------------------------------------------------
$dp = $wb.Worksheets.Item("Dynamic Programming")
$dp.Range("H5").Value = "The space complexity can be improved."
$dp.Range("H5").WrapText = $true
$dp.Range("E7").Value = "O(nm)"
$dp.Range("E7").WrapText = $true
$dp.Range("F7").Value = "O(nm)"
$dp.Range("F7").WrapText = $true
$dp.Columns.Item(8).ColumnWidth = 32.375

# ---------------------------------------------------------------------------
# 2. Strings sheet: add a new row for Caesar Cipher Encryptor
# ---------------------------------------------------------------------------
$strings = $wb.Worksheets.Item("Strings")
$strings.Range("A3").Value = "Caesar Cipher Encryptor"
$strings.Range("B3").Value = "Given a non-empty string of lowercase letters and a non-negative integer representing a key, write a function that returns a new string obtained by shifting every letter in the input string by k positions in the alphabet, where k is the key."
$strings.Range("B3").WrapText = $true
$strings.Range("C3").Value = "str = ""xyz"", k=2"
$strings.Range("D3").Value = """zab"""
$strings.Range("E3").Value = "O(n)"
$strings.Range("F3").Value = "O(n)"
$strings.Range("G3").Value = "ord('a') = 97 ; chr(97) = 'a'  list(""abc"") = ['a','b','c']"
$strings.Range("G3").WrapText = $true
$strings.Columns.Item(7).ColumnWidth = 30.75

# ---------------------------------------------------------------------------
# 3. New "Heaps" sheet (inserted after Strings)
# ---------------------------------------------------------------------------
$heaps = $wb.Worksheets.Add($null, $strings)
$heaps.Name = "Heaps"

$heaps.Range("A1").Value = "Problems"
$heaps.Range("A1").Font.Bold = $true
$heaps.Range("B1").Value = "Resources"
$heaps.Range("B1").Font.Bold = $true
$heaps.Range("B1").WrapText = $true
$heaps.Range("C1").Value = "Input"
$heaps.Range("C1").Font.Bold = $true
$heaps.Range("D1").Value = "Output"
$heaps.Range("D1").Font.Bold = $true
$heaps.Range("E1").Value = "Time"
$heaps.Range("E1").Font.Bold = $true
$heaps.Range("E1").WrapText = $true
$heaps.Range("F1").Value = "Space"
$heaps.Range("F1").Font.Bold = $true
$heaps.Range("F1").WrapText = $true
$heaps.Range("G1").Value = "Python Concepts"
$heaps.Range("G1").Font.Bold = $true
$heaps.Range("G1").WrapText = $true

$heaps.Range("A2").Value = "Min Heap"
$heaps.Range("E2").Value = "Peek O(1) BuildHeap - O(n)  SiftDown/SiftUp/Insert/Remove - O(log(n))  "
$heaps.Range("E2").WrapText = $true
$heaps.Range("F2").Value = "All - O(1)"
$heaps.Range("F2").WrapText = $true
$heaps.Range("G2").Value = " Decrementing For loop: for i in reversed(range(6)):"
$heaps.Range("G2").WrapText = $true
$heaps.Range("G2").Characters(2, 23).Font.Bold = $true

$heaps.Columns.Item(1).ColumnWidth = 23.625
$heaps.Columns.Item(2).ColumnWidth = 34.375
$heaps.Columns.Item(3).ColumnWidth = 19.5
$heaps.Columns.Item(4).ColumnWidth = 12
$heaps.Columns.Item(5).ColumnWidth = 23
$heaps.Columns.Item(6).ColumnWidth = 24
$heaps.Columns.Item(7).ColumnWidth = 34.875

$heaps.Range("B15").Select()

# ---------------------------------------------------------------------------
# 4. New "Graphs" sheet (inserted after Heaps)
# ---------------------------------------------------------------------------
$graphs = $wb.Worksheets.Add($null, $heaps)
$graphs.Name = "Graphs"

$graphs.Range("A1").Value = "Problems"
$graphs.Range("A1").Font.Bold = $true
$graphs.Range("B1").Value = "Resources"
$graphs.Range("B1").Font.Bold = $true
$graphs.Range("B1").WrapText = $true
$graphs.Range("C1").Value = "Input"
$graphs.Range("C1").Font.Bold = $true
$graphs.Range("C1").WrapText = $true
$graphs.Range("D1").Value = "Output"
$graphs.Range("D1").Font.Bold = $true
$graphs.Range("D1").WrapText = $true
$graphs.Range("E1").Value = "Time"
$graphs.Range("E1").Font.Bold = $true
$graphs.Range("F1").Value = "Space"
$graphs.Range("F1").Font.Bold = $true
$graphs.Range("F1").WrapText = $true
$graphs.Range("G1").Value = "Python Concepts"
$graphs.Range("G1").Font.Bold = $true

$graphs.Range("A2").Value = "Depth First Search"
$graphs.Range("E2").Value = "O(v+e)"
$graphs.Range("F2").Value = "O(v)"

$graphs.Columns.Item(1).ColumnWidth = 19.5
$graphs.Columns.Item(2).ColumnWidth = 18
$graphs.Columns.Item(3).ColumnWidth = 17.5
$graphs.Columns.Item(4).ColumnWidth = 13
$graphs.Columns.Item(5).ColumnWidth = 13.25
$graphs.Columns.Item(7).ColumnWidth = 22.625

$graphs.Range("A1:XFD1").Select()

# ---------------------------------------------------------------------------
# 5. Reorder "BinaryTrees" to the end (after Recursion)
# ---------------------------------------------------------------------------
$binaryTrees = $wb.Worksheets.Item("BinaryTrees")
$recursion = $wb.Worksheets.Item("Recursion")
$binaryTrees.Move($null, $recursion)

# ---------------------------------------------------------------------------
# 6. Selections on other, pre-existing sheets
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("D5").Select()

$array = $wb.Worksheets.Item("Array")
$array.Range("C6").Select()

$dp.Range("H6").Select()

$bst = $wb.Worksheets.Item("BinarySearchTrees")
$bst.Range("A1:XFD1").Select()

$searching = $wb.Worksheets.Item("Searching")
$searching.Range("A1:XFD1").Select()

$strings.Range("G5").Select()

$stacks = $wb.Worksheets.Item("Stacks")
$stacks.Range("B15").Select()

$recursion.Range("A3").Select()

$binaryTrees.Range("A1:XFD1").Select()

# ---------------------------------------------------------------------------
# 7. Active sheet / active cell for the whole workbook
# ---------------------------------------------------------------------------
$array.Activate()
$array.Range("C6").Select()
